# Apply the "cambios de agosto" edits to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: turn off wrap text on the merged description cell (G3:I3) and
#     bump the row height slightly to fit the single-line text ---
$ws.Range("G3:I3").WrapText = $false
$ws.Rows.Item(3).RowHeight = 34.5

# --- Row 8: shift the reporting period forward (Jan-Mar -> Apr-Jun 2022) ---
$ws.Range("B8").Value = 44652   # period start
$ws.Range("C8").Value = 44742   # period end
$ws.Range("E8").Value = 44727   # validation date

# F8 label: "Estadística enero-abril 2022" -> "Estadística mayo-agosto 2022"
$ws.Range("F8").Value = "Estadística mayo-agosto 2022"

$ws.Range("M8").Value = 44753   # fecha de validación
$ws.Range("N8").Value = 44753   # fecha de actualización

# --- Row 8 hyperlinks: only the J8 target URL actually changes (new zip
#     report for Abril-Junio), but this runtime's Hyperlink.Address setter
#     always appends a brand-new relationship/entry instead of editing one
#     in place (leaving the stale one behind), so stash the existing cell
#     formatting, rebuild the full hyperlink set in original order, then
#     restore the original per-cell styles that Hyperlinks.Add() disturbs. ---
$ws.Range("F8:L8").Copy() | Out-Null
$ws.Range("F30:L30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$urlG8 = "https://sistemas.dgutyp.sep.gob.mx/sii/Downloads/UTUPSII.pdf"
$urlH8 = "https://www.dof.gob.mx/nota_detalle.php?codigo=5600454&fecha=15/09/2020"
$urlK8 = "https://www.gob.mx/sep/acciones-y-programas/estadistica-educativa-15782"
$urlJ8 = "http://transparenciadocs.hidalgo.gob.mx/ENTIDADES/UPPachuca/dir1/2022/Abril-Junio/30/REPORTES%20SISTEMA%20MAYO%20AGOSTO%202022.zip"

# G8/H8/K8 keep their original display text (same as their Address, already
# present in the cell), so TextToDisplay is left out for them; only J8's
# visible text actually changes, so it is passed explicitly there.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G8"), $urlG8) | Out-Null
$ws.Hyperlinks.Add($ws.Range("H8"), $urlH8) | Out-Null
$ws.Hyperlinks.Add($ws.Range("K8"), $urlK8) | Out-Null
$ws.Hyperlinks.Add($ws.Range("J8"), $urlJ8, "", "", $urlJ8) | Out-Null

$ws.Range("F30:L30").Copy() | Out-Null
$ws.Range("F8:L8").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(30).Delete()

# --- Column widths ---
# This runtime quantizes ColumnWidth to 1/6-character steps and adds a fixed
# 5/6 offset when round-tripping through the OOXML `width` attribute
# (observed: saved_width = ceil(set_width*6)/6 + 5/6). Pre-compensate the
# value we set so the saved width lands on (or as close as achievable to)
# the target widths from the diff.
$ws.Columns.Item(4).ColumnWidth = 29.333333333333332    # -> width ~30.140625
$ws.Columns.Item(6).ColumnWidth = 27.666666666666668    # -> width ~28.42578125
$ws.Columns.Item(7).ColumnWidth = 66.16666666666667     # -> width 67 (exact)
$ws.Columns.Item(8).ColumnWidth = 81.5                  # -> width ~82.28515625
$ws.Columns.Item(10).ColumnWidth = 74.0                 # -> width ~74.85546875
$ws.Columns.Item(15).ColumnWidth = 15.166666666666666   # -> width 16 (exact)

# --- View state: scroll so column G is left-most, select H12 ---
$excel.ActiveWindow.ScrollColumn = 7   # G
$ws.Range("H12").Select()

$wb.Save()
